# Add a new date column (AF) "18-jul" with one day's worth of data,
# following the same pattern as the existing "17-jul" column (AE).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Best-effort: restore VBA-style code names (cosmetic metadata seen in the
# diff). Harmless if the host doesn't persist it.
try { $wb.CodeName = "ThisWorkbook" } catch {}
try { $ws.CodeName = "Sheet1" } catch {}

# New header cell for the added date column.
$ws.Range("AF1").Value = "18-jul"

# New data values for the added column (rows 2-18).
$ws.Range("AF2").Value  = 0
$ws.Range("AF3").Value  = 10.774470971714203
$ws.Range("AF4").Value  = 14.141731489157788
$ws.Range("AF5").Value  = 27.275380305916332
$ws.Range("AF6").Value  = 0
$ws.Range("AF7").Value  = 4.2103534177627351
$ws.Range("AF8").Value  = 13.535980989973574
$ws.Range("AF9").Value  = 22.358540510919667
$ws.Range("AF10").Value = 21.852220645540704
$ws.Range("AF11").Value = 9.6534038810884795
$ws.Range("AF12").Value = 0
$ws.Range("AF13").Value = 7.8441604470884752
$ws.Range("AF14").Value = 0
$ws.Range("AF15").Value = 0
$ws.Range("AF16").Value = 15.817764657367206
$ws.Range("AF17").Value = 0
$ws.Range("AF18").Value = 0

# Match the updated selection recorded in the workbook.
$ws.Range("AB8").Select()
